$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster label columns reference (shared strings): ECs, FAPs, M2, sCs
$clusters = @("ECs", "FAPs", "M2", "sCs")
$targets  = @("ECs", "M2")

# Data rows 2-9: sending cluster x target cluster grid
$data = @(
    # A(sending), D(target), E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
    @("ECs","ECs",3,1,21.557693,64.673079,0.5505707555812251,0.5505707555812251,1,0.3333333333333333,0.1791996666666667,0.537599,0.04251079199666429,0.04251079199666429,3.863131399702334,34.768182597321,0.02340519886995976,0.02340519886995976),
    @("ECs","M2",3,1,21.557693,64.673079,0.5505707555812251,0.5505707555812251,3,1,4.036192666666667,12.108578,0.9574892080033358,0.9574892080033357,87.01100239685134,783.099021571662,0.5271655567112654,0.5271655567112654),
    @("FAPs","ECs",3,1,3.981869333333333,11.945608,0.1016945926207894,0.1016945926207894,1,0.3333333333333333,0.1791996666666667,0.537599,0.04251079199666429,0.04251079199666429,0.7135496572435556,6.421946915192001,0.004323117674087888,0.004323117674087888),
    @("FAPs","M2",3,1,3.981869333333333,11.945608,0.1016945926207894,0.1016945926207894,3,1,4.036192666666667,12.108578,0.9574892080033358,0.9574892080033357,16.07159180282489,144.644326225424,0.09737147494670148,0.09737147494670148),
    @("M2","ECs",3,1,12.00696933333334,36.02090800000001,0.3066509100994217,0.3066509100994217,1,0.3333333333333333,0.1791996666666667,0.537599,0.04251079199666429,0.04251079199666429,2.151644902210223,19.364804119892,0.01303597305482432,0.01303597305482432),
    @("M2","M2",3,1,12.00696933333334,36.02090800000001,0.3066509100994217,0.3066509100994217,3,1,4.036192666666667,12.108578,0.9574892080033358,0.9574892080033357,48.46244157209156,436.1619741488241,0.2936149370445975,0.2936149370445974),
    @("sCs","ECs",3,1,1.608641,4.825923,0.04108374169856382,0.04108374169856382,1,0.3333333333333333,0.1791996666666667,0.537599,0.04251079199666429,0.04251079199666429,0.2882679309863333,2.594411378877,0.00174650239779233,0.00174650239779233),
    @("sCs","M2",3,1,1.608641,4.825923,0.04108374169856382,0.04108374169856382,3,1,4.036192666666667,12.108578,0.9574892080033358,0.9574892080033357,6.492785007499332,58.43506506749399,0.0393372393007715,0.03933723930077149)
)

$rowIdx = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIdx, 1).Value = $rec[0]
    $ws.Cells.Item($rowIdx, 2).Value = "Cd80"
    $ws.Cells.Item($rowIdx, 3).Value = "Cd28"
    $ws.Cells.Item($rowIdx, 4).Value = $rec[1]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($rowIdx, $c).Value = $rec[$c - 3]
    }
    $rowIdx++
}
